$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2025-10-28 Tuesday" "2025-10-29 Wednesday"

Replace-Text "542÷2=" "103÷7="
Replace-Text "101÷2=" "401÷9="
Replace-Text "627÷5=" "454÷6="
Replace-Text "321÷9=" "219÷4="
Replace-Text "407÷6=" "539÷8="
Replace-Text "691÷5=" "588÷9="
Replace-Text "912÷7=" "991÷4="
Replace-Text "787÷7=" "892÷7="
Replace-Text "830÷6=" "702÷5="
Replace-Text "341÷6=" "712÷8="
Replace-Text "558÷8=" "260÷3="
Replace-Text "445÷3=" "377÷2="
Replace-Text "158÷9=" "299÷5="
Replace-Text "862÷2=" "867÷4="
Replace-Text "501÷8=" "790÷2="
Replace-Text "608÷3=" "296÷2="
Replace-Text "147÷5=" "658÷8="
Replace-Text "115÷6=" "963÷6="
Replace-Text "652÷8=" "770÷7="
Replace-Text "725÷4=" "869÷4="
Replace-Text "682÷2=" "684÷4="
Replace-Text "617÷7=" "692÷4="
Replace-Text "826÷8=" "664÷6="
Replace-Text "799÷6=" "713÷4="
Replace-Text "267÷2=" "934÷6="
